$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update to new values (columns B..V), I2 flips from 0 to 1
$ws.Range("B2:H2").Value = 0
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0.06313584232034301
$ws.Range("K2").Value = 0.0485201827700498
$ws.Range("L2").Value = 0.03282828711956988
$ws.Range("N2").Value = 1037.038442353979
$ws.Range("P2").Value = 0.9051765380291221
$ws.Range("R2").Value = 0.8367506100813272
$ws.Range("T2").Value = 0.3697803381425155
$ws.Range("V2").Value = 0.4748804363443177

# Rows 3..11: zero out every value column (B..V); columns already 0 are unaffected
$ws.Range("B3:V11").Value = 0
